$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.84"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "11"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.55"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "11"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.381"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "11"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05762"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "11"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "11"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.337"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "11"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "11"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8946"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "11"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1446"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "11"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07348"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11BitrueCoinBTRBestin24h"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "11"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02972"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "11"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09409"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "11"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.950"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "11"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001601"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "11"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04798"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "11"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "One"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005843"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "11"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006168"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "11"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004064"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "11"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009929"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "11"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "11"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.746"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "11"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.196"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "11"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3274"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "11"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "11"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003160"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "11"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "11"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "11"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "11"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "11"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "11"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "11"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "11"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "11"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "11"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "11"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "11"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "11"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03895"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "11"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006792"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "11"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1072"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "11"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002411"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "11"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006931"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "11"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005641"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "11"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "11"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3802"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "11"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1657"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "11"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "11"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "11"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "11"

Write-Output "Applied 118 cell updates"